$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Duplicate Sheet1 into a new sheet positioned right after it, named "Sheet_2"
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet_2"

# Update the scenario-specific inputs on the new sheet
$ws2.Range("K3").Value = 30000
$ws2.Range("J4").Value = 30000
$ws2.Range("D6").Value = 0.235
$ws2.Range("D7").Formula = '=(($K$3)*$L$3+(D4-$J$4)*$L$4)/D4'
$ws2.Range("D7").NumberFormat = $ws2.Range("E7").NumberFormat

# Sheet1: give L3 the percent number format used by its neighbour L4
$ws1.Range("L3").NumberFormat = $ws1.Range("L4").NumberFormat

# Restore the cursor positions recorded for each sheet, and make Sheet_2 the active tab
[void]$ws1.Range("E12").Select()
[void]$ws2.Range("G14").Select()
$ws2.Activate()
